$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Add new row 64: "1926. Nearest Exit from Entrance in Maze" ---
$ws.Range("A64").Value = "1926. Nearest Exit from Entrance in Maze"
$ws.Range("B64").Value = "Medium"
$ws.Range("C64").Value = "Graphs"
$ws.Range("D64").Value = "Classic Graph BFS. Construct the adjacency list, then perform BFS and check conditions. Return the level when a valid goal node is found."

$linkText = "https://leetcode.com/problems/nearest-exit-from-entrance-in-maze/solutions/2834640/java-explained-in-detail-simple-fast-solution-bfs/?envType=study-plan-v2&envId=leetcode-75 "
$linkAddress = $linkText.Trim()
$ws.Range("E64").Value = $linkText

# Wire up the hyperlink (mirrors the other rows in column E)
$ws.Hyperlinks.Add($ws.Range("E64"), $linkAddress, "", "", $linkText) | Out-Null

# Re-apply the same cell formatting used by the preceding data row so the
# "Medium" fill (B) and Hyperlink style (E) match the rest of the table
# (Hyperlinks.Add nudges the xf it assigns, so fix it up afterwards).
$ws.Range("B63").Copy()
$ws.Range("B64").PasteSpecial(-4122)
$ws.Range("E63").Copy()
$ws.Range("E64").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update the view/selection state to match the edited workbook ---
# (scrolled one column to the right: topLeftCell A40 -> B40, selection D67 -> E67)
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 40
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("E67").Select()
